# Updated latest Guinea data - loc_holiday,Zone,Center.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row (row 11) with a new Guinea holiday entry.
$row = 11

$ws.Cells.Item($row, 1).Value = 2000010
$ws.Cells.Item($row, 2).Value = "GN"
$ws.Cells.Item($row, 3).Value = 45567
$ws.Cells.Item($row, 4).Value = "Guinea day"
$ws.Cells.Item($row, 5).Value = "Guinea day"
$ws.Cells.Item($row, 6).Value = "fra"
$ws.Cells.Item($row, 7).Value = $true
$ws.Cells.Item($row, 8).Value = "superadmin"
$ws.Cells.Item($row, 9).Value = 45224.547017106481
$ws.Cells.Item($row, 10).Value = "NULL"
$ws.Cells.Item($row, 11).Value = "NULL"
$ws.Cells.Item($row, 12).Value = $false
$ws.Cells.Item($row, 13).Value = "NULL"

# Reuse the existing number-format styles from the row above instead of
# creating brand-new style entries.
$ws.Cells.Item($row - 1, 3).Copy()
$ws.Cells.Item($row, 3).PasteSpecial(-4122)

$ws.Cells.Item($row - 1, 9).Copy()
$ws.Cells.Item($row, 9).PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update sheet view to match target: topLeftCell H1, selection M14
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("M14").Select()
